$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tab-separated-like (pipe-delimited) table of target values for rows 2-51:
# RowNumber|Coin|Link|Price|Volume(1h)
$data = @"
2|Bitcoin|https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc|75.750.69|  +0.55%  
3|Ethereum|https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth|2.905.03|  +2.69%  
4|TetherUSD|https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt|1.00|  +0.05%  
5|Solana|https://coinranking.com/coin/zNZHO_Sjf+solana-sol|197.01|  +3.26%  
6|BNB|https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb|597.37|  -1.26%  
7|USDC|https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc|0.999|  -0.02%  
8|XRP|https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp|0.551|  +0.98%  
9|Dogecoin|https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge|0.192|  -4.10%  
10|LidoStakedEther|https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth|2.903.75|  +2.84%  
11|Cardano|https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada|0.413|  +10.38%  
12|TRON|https://coinranking.com/coin/qUhEFk1I61atv+tron-trx|0.160|  -1.21%  
13|Toncoin|https://coinranking.com/coin/67YlI0K1b+toncoin-ton|4.89|  +0.47%  
14|WrappedliquidstakedEther2.0|https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth|3.440.72|  +2.78%  
15|WrappedBTC|https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc|75.718.91|  +0.61%  
16|Avalanche|https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax|27.36|  -1.57%  
17|ShibaInu|https://coinranking.com/coin/xz24e0BjL+shibainu-shib|0.0000189|  -2.50%  
18|WrappedEther|https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth|2.913.19|  +2.96%  
19|Uniswap|https://coinranking.com/coin/_H5FVG9iW+uniswap-uni|8.89|  -5.46%  
20|Chainlink|https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link|12.62|  +0.70%  
21|BitcoinCash|https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch|376.26|  -1.90%  
22|SuiNetwork|https://coinranking.com/coin/3xJluUMvp+suinetwork-sui|2.30|  -0.74%  
23|Polkadot|https://coinranking.com/coin/25W7FG7om+polkadot-dot|4.14|  -1.16%  
24|Litecoin|https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc|71.40|  +0.08%  
25|Dai|https://coinranking.com/coin/MoTuySvg7+dai-dai|1.00|  -0.08%  
26|WrappedeETH|https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth|3.063.37|  +3.47%  
27|NEARProtocol|https://coinranking.com/coin/DCrsaMv68+nearprotocol-near|4.22|  -2.42%  
28|Aptos|https://coinranking.com/coin/HGYj5JCv5+aptos-apt|9.62|  -1.64%  
29|PEPE|https://coinranking.com/coin/03WI8NQPF+pepe-pepe|0.0000109|  +2.87%  
30|Binance-PegBSC-USD|https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd|1.00|  +0.21%  
31|Fetch.AI|https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet|1.40|  -2.96%  
32|Bittensor|https://coinranking.com/coin/pgv7xSFi6+bittensor-tao|502.75|  -6.26%  
33|InternetComputer(DFINITY)|https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp|7.72|  -4.08%  
34|PancakeSwap|https://coinranking.com/coin/ncYFcP709+pancakeswap-cake|1.81|  -2.06%  
35|FirstDigitalUSD|https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd|1.00|  +0.12%  
36|Monero|https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr|163.14|  +0.20%  
37|EthereumClassic|https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc|20.01|  -2.04%  
38|WhiteBITCoin|https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt|19.69|  +1.88%  
39|Kaspa|https://coinranking.com/coin/V8GxkwWow+kaspa-kas|0.113|  -7.26%  
40|USDe|https://coinranking.com/coin/exbfr2U-0+usde-usde|1.00|  -0.06%  
41|Aave|https://coinranking.com/coin/ixgUfzmLR+aave-aave|180.40|  -3.74%  
42|PolygonEcosystemToken|https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol|0.344|  -0.11%  
43|RenderToken|https://coinranking.com/coin/vfo5XYwcV+rendertoken-render|5.00|  -4.21%  
44|Stacks|https://coinranking.com/coin/mMPrMcB7+stacks-stx|1.66|  -4.47%  
45|Cronos|https://coinranking.com/coin/65PHZTpmE55b+cronos-cro|0.0907|  +5.31%  
46|ImmutableX|https://coinranking.com/coin/Z96jIvLU7+immutablex-imx|1.21|  -4.74%  
47|OKB|https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb|40.03|  +0.26%  
48|dogwifhat|https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif|2.39|  -2.11%  
49|ARBITRUM|https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb|0.577|  -0.66%  
50|Filecoin|https://coinranking.com/coin/ymQub4fuB+filecoin-fil|3.72|  -2.35%  
51|Mantle|https://coinranking.com/coin/BoI4ux0nd+mantle-mnt|0.654|  +7.06%  
"@

$lines = $data -split "`r`n|`n"

foreach ($line in $lines) {
    if ($line.Length -eq 0) { continue }
    $parts = $line.Split("|")
    $rowNum = [int]$parts[0]
    $coin = $parts[1]
    $link = $parts[2]
    $price = $parts[3]
    $volume = $parts[4]

    # Column B: Coin name (text)
    $ws.Cells.Item($rowNum, 2).Value = $coin

    # Column C: Link (text)
    $ws.Cells.Item($rowNum, 3).Value = $link

    # Column D: Price - force as text so values like "1.00" or "75.750.69"
    # are not reinterpreted as numbers, then restore default styling.
    $priceCell = $ws.Cells.Item($rowNum, 4)
    $origStyle = $priceCell.Style
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.Style = $origStyle

    # Column E: Volume(1h) - also force as text to preserve formatting
    # (leading/trailing spaces, % sign).
    $volCell = $ws.Cells.Item($rowNum, 5)
    $origStyleE = $volCell.Style
    $volCell.NumberFormat = "@"
    $volCell.Value = $volume
    $volCell.Style = $origStyleE
}
